# Updated cryptos list data (prices + 1h volume change) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '67.453.58'
$ws.Range("E2").Value = '  -0.67%  '

$ws.Range("D3").Value = "'" + '3.518.69'
$ws.Range("E3").Value = '  -1.24%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'" + '610.80'
$ws.Range("E5").Value = '  -1.09%  '

$ws.Range("D6").Value = "'" + '150.64'
$ws.Range("E6").Value = '  -2.14%  '

$ws.Range("D7").Value = "'" + '3.519.12'
$ws.Range("E7").Value = '  -1.06%  '

$ws.Range("D9").Value = "'" + '0.482'
$ws.Range("E9").Value = '  -1.20%  '

$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("E11").Value = '  +1.83%  '

$ws.Range("D12").Value = "'" + '0.426'
$ws.Range("E12").Value = '  -1.55%  '

$ws.Range("D13").Value = "'" + '0.0000220'
$ws.Range("E13").Value = '  -2.15%  '

$ws.Range("D14").Value = "'" + '4.114.34'
$ws.Range("E14").Value = '  -1.25%  '

$ws.Range("D15").Value = "'" + '31.81'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = "'" + '3.507.11'
$ws.Range("E16").Value = '  -1.64%  '

$ws.Range("D17").Value = "'" + '67.433.77'
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").Value = "'" + '15.24'
$ws.Range("E20").Value = '  -2.59%  '

$ws.Range("D21").Value = "'" + '443.40'
$ws.Range("E21").Value = '  -2.81%  '

$ws.Range("D22").Value = "'" + '9.27'
$ws.Range("E22").Value = '  -3.99%  '

$ws.Range("E23").Value = '  -3.38%  '

$ws.Range("D24").Value = "'" + '77.26'
$ws.Range("E24").Value = '  -0.49%  '

$ws.Range("D25").Value = "'" + '0.0000129'
$ws.Range("E25").Value = '  +11.40%  '

$ws.Range("D26").Value = "'" + '3.659.31'
$ws.Range("E26").Value = '  -1.23%  '

$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").Value = "'" + '10.19'
$ws.Range("E28").Value = '  -4.33%  '

$ws.Range("D29").Value = "'" + '8.35'
$ws.Range("E29").Value = '  +0.22%  '

$ws.Range("E30").Value = '  -2.30%  '

$ws.Range("E31").Value = '  -4.85%  '

$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("D33").Value = "'" + '0.164'
$ws.Range("E33").Value = '  +3.49%  '

$ws.Range("D34").Value = "'" + '25.79'
$ws.Range("E34").Value = '  -0.78%  '

$ws.Range("D35").Value = "'" + '6.14'
$ws.Range("E35").Value = '  -1.12%  '

$ws.Range("D36").Value = "'" + '3.511.05'
$ws.Range("E36").Value = '  -1.46%  '

$ws.Range("E37").Value = '  -3.50%  '

$ws.Range("D38").Value = "'" + '8.04'
$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = "'" + '178.08'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = "'" + '1.00'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").Value = "'" + '2.18'
$ws.Range("E42").Value = '  +4.05%  '

$ws.Range("D43").Value = "'" + '0.0874'
$ws.Range("E43").Value = '  -1.24%  '

$ws.Range("E44").Value = '  -3.23%  '

$ws.Range("D45").Value = "'" + '0.881'
$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("D46").Value = "'" + '45.57'
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("D47").Value = "'" + '27.71'
$ws.Range("E47").Value = '  -2.49%  '

$ws.Range("D48").Value = "'" + '1.26'
$ws.Range("E48").Value = '  +4.43%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("D51").Value = "'" + '0.995'
$ws.Range("E51").Value = '  -1.62%  '
